# Regenerate sval data to filter save games.
# Updates the computed stat columns (B-E, G) for rows 2-8 on the active sheet
# with the refreshed values from the new data pipeline run. Column F (Win)
# is left untouched since it is unaffected by this regeneration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.01253208636536152; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 2.89400026249618 }
    3 = @{ B = 1.445647641019636;   C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
    4 = @{ B = 3.272327238179451;   C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
    5 = @{ B = 3.272327238179451;   C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
    6 = @{ B = 3.272327238179451;   C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
    7 = @{ B = 0.6545652718822623;  C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 6.038307959104277 }
    8 = @{ B = 0.6545652718822623;  C = 0.04103571897497393; D = 0.1496068669990043; E = 0.5333859586016987; G = 1.378593816457939 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
